# Applies the "Updated ITA model" edit to scen_tsparameters_ts48_clu.xlsx
#  - ev_charging_uc!C13 / C14: re-shuffled (same set, new order) comma-separated
#    timeslice lists used by the Day / Night UC_T rows
#  - re_profiles!N11:N50: previously-zero shares filled in with computed values
#  - re_profiles!Q13:R14: the S3 / S4 season rows swapped places

$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": refresh the Day (C13) / Night (C14) timeslice lists ---
$ws1 = $wb.Worksheets.Item("ev_charging_uc")
$ws1.Range("C13").Value = "S4aH2,S1aH3,S3aH2,S3aH3,S5aH4,S3aH7,S4aH5,S1aH6,S2aH6,S3aH5,S5aH3,S2aH7,S3aH4,S4aH3,S4aH4,S1aH2,S1aH7,S2aH3,S3aH6,S4aH7,S5aH6,S5aH7,S1aH4,S2aH2,S4aH6,S1aH5,S2aH4,S2aH5,S5aH2,S5aH5"
$ws1.Range("C14").Value = "S2aH1,S5aH8,S1aH1,S3aH8,S4aH8,S3aH1,S2aH8,S1aH8,S5aH1,S4aH1"

# --- Sheet "re_profiles": fill in column N (rows 11-50) ---
$ws4 = $wb.Worksheets.Item("re_profiles")
$ws4.Cells.Item(11, 14).Value = 0.0682325241216655
$ws4.Cells.Item(12, 14).Value = 0.010064341039034391
$ws4.Cells.Item(13, 14).Value = 0.010053840550443846
$ws4.Cells.Item(14, 14).Value = 0.0099879348718522783
$ws4.Cells.Item(15, 14).Value = 0.060792920303432134
$ws4.Cells.Item(16, 14).Value = 0.010040614148058175
$ws4.Cells.Item(17, 14).Value = 0.0099745987670680752
$ws4.Cells.Item(18, 14).Value = 0.058119770666056537
$ws4.Cells.Item(19, 14).Value = 0.038375118931494201
$ws4.Cells.Item(20, 14).Value = 0.0053450624963613629
$ws4.Cells.Item(21, 14).Value = 0.0052757724208315829
$ws4.Cells.Item(22, 14).Value = 0.0052449034015103913
$ws4.Cells.Item(23, 14).Value = 0.032253569290825734
$ws4.Cells.Item(24, 14).Value = 0.0055064871848617535
$ws4.Cells.Item(25, 14).Value = 0.0055870317896012491
$ws4.Cells.Item(26, 14).Value = 0.032904049220948452
$ws4.Cells.Item(27, 14).Value = 0.13470775990084147
$ws4.Cells.Item(28, 14).Value = 0.017906110932828567
$ws4.Cells.Item(29, 14).Value = 0.017662605018771465
$ws4.Cells.Item(30, 14).Value = 0.017541144373550716
$ws4.Cells.Item(31, 14).Value = 0.11075349282913861
$ws4.Cells.Item(32, 14).Value = 0.020423104643675052
$ws4.Cells.Item(33, 14).Value = 0.020688393709751108
$ws4.Cells.Item(34, 14).Value = 0.12284491767933799
$ws4.Cells.Item(35, 14).Value = 0.03076516287723392
$ws4.Cells.Item(36, 14).Value = 0.0041971271016678496
$ws4.Cells.Item(37, 14).Value = 0.0041004062264624781
$ws4.Cells.Item(38, 14).Value = 0.0040255399778688682
$ws4.Cells.Item(39, 14).Value = 0.023322213983746443
$ws4.Cells.Item(40, 14).Value = 0.0040592208003663978
$ws4.Cells.Item(41, 14).Value = 0.0041210425341255304
$ws4.Cells.Item(42, 14).Value = 0.026532671361551369
$ws4.Cells.Item(43, 14).Value = 0.02048320543136058
$ws4.Cells.Item(44, 14).Value = 0.0027219450284947876
$ws4.Cells.Item(45, 14).Value = 0.0026982214942662254
$ws4.Cells.Item(46, 14).Value = 0.0027104471926864219
$ws4.Cells.Item(47, 14).Value = 0.016573906647167156
$ws4.Cells.Item(48, 14).Value = 0.0028729049965876639
$ws4.Cells.Item(49, 14).Value = 0.0029121969268978559
$ws4.Cells.Item(50, 14).Value = 0.01761771912737203

# Swap the S3 / S4 rows (Q/R columns) on rows 13-14
$ws4.Cells.Item(13, 17).Value = "S4"
$ws4.Cells.Item(13, 18).Value = 0.084821841739935219
$ws4.Cells.Item(14, 17).Value = "S3"
$ws4.Cells.Item(14, 18).Value = 0.81410226746876446
